$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new header labels in row 1
$ws.Range("K1").Value = "svr_kernel_scale"
$ws.Range("L1").Value = "svr_epsilon"
$ws.Range("M1").Value = "svr_box_constraint"

# Add new data values in row 2
$ws.Range("K2").Value = 100
$ws.Range("L2").Value = 0.1
$ws.Range("M2").Value = 5

# Update selection to match final state
$ws.Range("L8").Select()
